# Updated cryptos list - applies latest price/volume scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT (matches source data, which stores
# prices/percentages as inline strings, not numbers) while leaving the
# cell's style untouched (reset back to "Normal" after the temporary
# "@" text format so no stray style survives the round-trip).
function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue $ws "D2" '87.782.41'
$ws.Range("E2").Value = '  +0.09%  '
Set-TextValue $ws "D3" '3.170.57'
$ws.Range("E3").Value = '  -4.10%  '
$ws.Range("E4").Value = '  -0.13%  '
Set-TextValue $ws "D5" '207.48'
$ws.Range("E5").Value = '  -3.31%  '
Set-TextValue $ws "D6" '610.49'
$ws.Range("E6").Value = '  -3.21%  '
$ws.Range("E7").Value = '  -2.39%  '
$ws.Range("E8").Value = '  +4.09%  '
$ws.Range("E9").Value = '  -0.10%  '
Set-TextValue $ws "D10" '3.167.61'
$ws.Range("E10").Value = '  -4.10%  '
Set-TextValue $ws "D11" '0.537'
$ws.Range("E11").Value = '  -9.62%  '
$ws.Range("E12").Value = '  -0.80%  '
$ws.Range("E13").Value = '  -8.65%  '
Set-TextValue $ws "D14" '3.760.18'
$ws.Range("E14").Value = '  -4.00%  '
Set-TextValue $ws "D15" '5.28'
$ws.Range("E15").Value = '  -0.43%  '
Set-TextValue $ws "D16" '87.789.89'
$ws.Range("E16").Value = '  +0.35%  '
Set-TextValue $ws "D17" '32.21'
$ws.Range("E17").Value = '  -7.49%  '
Set-TextValue $ws "D18" '3.161.22'
$ws.Range("E18").Value = '  -4.26%  '
Set-TextValue $ws "D19" '3.18'
$ws.Range("E19").Value = '  +5.14%  '
Set-TextValue $ws "D20" '13.47'
$ws.Range("E20").Value = '  -5.92%  '
Set-TextValue $ws "D21" '413.26'
$ws.Range("E21").Value = '  -5.92%  '
Set-TextValue $ws "D22" '8.47'
$ws.Range("E22").Value = '  -8.52%  '
Set-TextValue $ws "D23" '5.06'
$ws.Range("E23").Value = '  -6.61%  '
Set-TextValue $ws "D24" '5.23'
$ws.Range("E24").Value = '  -0.64%  '
Set-TextValue $ws "D25" '12.19'
$ws.Range("E25").Value = '  -2.01%  '
Set-TextValue $ws "D26" '3.344.66'
$ws.Range("E26").Value = '  -3.57%  '
Set-TextValue $ws "D27" '0.0000132'
$ws.Range("E27").Value = '  -0.19%  '
Set-TextValue $ws "D28" '73.47'
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("E30").Value = '  -12.33%  '
Set-TextValue $ws "D31" '1.00'
$ws.Range("E31").Value = '  -0.79%  '
Set-TextValue $ws "D32" '545.66'
$ws.Range("E32").Value = '  -2.12%  '
Set-TextValue $ws "D33" '8.23'
$ws.Range("E33").Value = '  -9.29%  '
Set-TextValue $ws "D34" '1.32'
$ws.Range("E34").Value = '  -9.68%  '
$ws.Range("E35").Value = '  -1.61%  '
$ws.Range("E36").Value = '  -6.78%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws "D37" '0.130'
$ws.Range("E37").Value = '  -7.36%  '
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws "D38" '21.90'
$ws.Range("E38").Value = '  -4.20%  '
Set-TextValue $ws "D39" '21.84'
$ws.Range("E39").Value = '  +0.29%  '
$ws.Range("E40").Value = '  -0.09%  '
Set-TextValue $ws "D41" '3.03'
$ws.Range("E41").Value = '  +0.96%  '
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("E43").Value = '  -6.99%  '
Set-TextValue $ws "D44" '0.373'
$ws.Range("E44").Value = '  -7.72%  '
Set-TextValue $ws "D45" '148.38'
$ws.Range("E45").Value = '  -4.18%  '
Set-TextValue $ws "D46" '173.98'
$ws.Range("E46").Value = '  -4.61%  '
Set-TextValue $ws "D47" '43.22'
$ws.Range("E47").Value = '  -4.64%  '
$ws.Range("E48").Value = '  +3.48%  '
$ws.Range("E49").Value = '  -10.05%  '
Set-TextValue $ws "D50" '3.96'
$ws.Range("E50").Value = '  -8.00%  '
Set-TextValue $ws "D51" '23.88'
$ws.Range("E51").Value = '  -3.59%  '
